$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 15
$ws.Range("BF2").Value = "'2013-01-31"
$ws.Range("AD3").Value = 15
$ws.Range("BF3").Value = "'2013-01-31"
$ws.Range("AH4").Value = 10
$ws.Range("BF4").Value = "'2013-01-31"
$ws.Range("AD5").Value = 15
$ws.Range("BF5").Value = "'2013-01-31"
$ws.Range("AD6").Value = 15
$ws.Range("AE6").Value = 7
$ws.Range("AG6").Value = 7
$ws.Range("AH6").Value = 14
$ws.Range("BF6").Value = "'2013-01-31"
$ws.Range("AW7").Value = 13
$ws.Range("BF7").Value = "'2013-01-31"
$ws.Range("D8").Value = 45
$ws.Range("F8").Value = 26
$ws.Range("G8").Value = 0.422
$ws.Range("I8").Value = 38.2
$ws.Range("J8").Value = 84.3
$ws.Range("M8").Value = 19.3
$ws.Range("N8").Value = 0.362
$ws.Range("P8").Value = 21.6
$ws.Range("Q8").Value = 0.798
$ws.Range("S8").Value = 32.5
$ws.Range("T8").Value = 42
$ws.Range("AA8").Value = 19.2
$ws.Range("AB8").Value = 100.6
$ws.Range("AD8").Value = 15
$ws.Range("AF8").Value = 19
$ws.Range("AG8").Value = 19
$ws.Range("AZ8").Value = 25
$ws.Range("BF8").Value = "'2013-01-31"
$ws.Range("AR9").Value = 1
$ws.Range("BF9").Value = "'2013-01-31"
$ws.Range("AH10").Value = 10
$ws.Range("BF10").Value = "'2013-01-31"
$ws.Range("D11").Value = 45
$ws.Range("E11").Value = 28
$ws.Range("G11").Value = 0.622
$ws.Range("K11").Value = 0.458
$ws.Range("N11").Value = 0.394
$ws.Range("O11").Value = 17
$ws.Range("P11").Value = 21.1
$ws.Range("Q11").Value = 0.804
$ws.Range("R11").Value = 10.8
$ws.Range("T11").Value = 44.6
$ws.Range("V11").Value = 15.3
$ws.Range("AD11").Value = 15
$ws.Range("AE11").Value = 7
$ws.Range("AH11").Value = 14
$ws.Range("AK11").Value = 8
$ws.Range("AO11").Value = 14
$ws.Range("AU11").Value = 8
$ws.Range("AV11").Value = 24
$ws.Range("BC11").Value = 12
$ws.Range("BF11").Value = "'2013-01-31"
$ws.Range("AT12").Value = 10
$ws.Range("BF12").Value = "'2013-01-31"
$ws.Range("AH13").Value = 16
$ws.Range("AV13").Value = 23
$ws.Range("BF13").Value = "'2013-01-31"
$ws.Range("AE14").Value = 2
$ws.Range("AO14").Value = 12
$ws.Range("BF14").Value = "'2013-01-31"
$ws.Range("AK15").Value = 9
$ws.Range("BF15").Value = "'2013-01-31"
$ws.Range("D16").Value = 44
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 0.659
$ws.Range("H16").Value = 48.5
$ws.Range("J16").Value = 82.8
$ws.Range("K16").Value = 0.437
$ws.Range("M16").Value = 13.7
$ws.Range("N16").Value = 0.344
$ws.Range("O16").Value = 16.3
$ws.Range("P16").Value = 20.6
$ws.Range("R16").Value = 13.6
$ws.Range("S16").Value = 29.4
$ws.Range("T16").Value = 43
$ws.Range("U16").Value = 20.6
$ws.Range("W16").Value = 9
$ws.Range("Y16").Value = 6
$ws.Range("AA16").Value = 19.5
$ws.Range("AB16").Value = 93.40000000000001
$ws.Range("AC16").Value = 3.9
$ws.Range("AD16").Value = 25
$ws.Range("AF16").Value = 5
$ws.Range("AG16").Value = 5
$ws.Range("AJ16").Value = 10
$ws.Range("AN16").Value = 24
$ws.Range("AR16").Value = 2
$ws.Range("AT16").Value = 12
$ws.Range("AV16").Value = 16
$ws.Range("AY16").Value = 22
$ws.Range("BB16").Value = 27
$ws.Range("BF16").Value = "'2013-01-31"
$ws.Range("BF17").Value = "'2013-01-31"
$ws.Range("AD18").Value = 25
$ws.Range("AI18").Value = 9
$ws.Range("BF18").Value = "'2013-01-31"
$ws.Range("AU19").Value = 18
$ws.Range("AY19").Value = 24
$ws.Range("BF19").Value = "'2013-01-31"
$ws.Range("AH20").Value = 16
$ws.Range("AY20").Value = 25
$ws.Range("BF20").Value = "'2013-01-31"
$ws.Range("AE21").Value = 7
$ws.Range("AG21").Value = 6
$ws.Range("AJ21").Value = 8
$ws.Range("BF21").Value = "'2013-01-31"
$ws.Range("D22").Value = 45
$ws.Range("E22").Value = 34
$ws.Range("G22").Value = 0.756
$ws.Range("H22").Value = 48.6
$ws.Range("I22").Value = 37.6
$ws.Range("J22").Value = 78.90000000000001
$ws.Range("K22").Value = 0.477
$ws.Range("L22").Value = 7.6
$ws.Range("M22").Value = 19.6
$ws.Range("P22").Value = 27.3
$ws.Range("R22").Value = 10.6
$ws.Range("S22").Value = 32.4
$ws.Range("X22").Value = 6.9
$ws.Range("Y22").Value = 4.1
$ws.Range("AC22").Value = 8.4
$ws.Range("AD22").Value = 15
$ws.Range("AH22").Value = 8
$ws.Range("AI22").Value = 10
$ws.Range("AM22").Value = 13
$ws.Range("AS22").Value = 7
$ws.Range("AT22").Value = 10
$ws.Range("AU22").Value = 19
$ws.Range("AW22").Value = 12
$ws.Range("BC22").Value = 2
$ws.Range("BF22").Value = "'2013-01-31"
$ws.Range("AD23").Value = 15
$ws.Range("AM23").Value = 14
$ws.Range("BF23").Value = "'2013-01-31"
$ws.Range("AD24").Value = 15
$ws.Range("BB24").Value = 28
$ws.Range("BF24").Value = "'2013-01-31"
$ws.Range("BF25").Value = "'2013-01-31"
$ws.Range("AD26").Value = 15
$ws.Range("BF26").Value = "'2013-01-31"
$ws.Range("AJ27").Value = 7
$ws.Range("AZ27").Value = 24
$ws.Range("BF27").Value = "'2013-01-31"
$ws.Range("BC28").Value = 1
$ws.Range("BF28").Value = "'2013-01-31"
$ws.Range("AO29").Value = 13
$ws.Range("BF29").Value = "'2013-01-31"
$ws.Range("AH30").Value = 16
$ws.Range("AU30").Value = 7
$ws.Range("AY30").Value = 23
$ws.Range("BF30").Value = "'2013-01-31"
$ws.Range("AD31").Value = 25
$ws.Range("AN31").Value = 25
$ws.Range("BF31").Value = "'2013-01-31"
